$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (45203) for every data
# row from row 2 through row 387. Bump every value from 45203 to 45204.
$ws.Range("C2:C387").Value2 = 45204
